# Apply the "topics_comparison" update: five new topic rows (7-11) get a
# label in column C and are filled in with "x" marks (re-using the existing
# blue "x" cell style) in the columns where that topic applies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New topic names for rows 9-13 (No. 7-11)
$ws.Range("C9").Value  = "Non-Destructive Testing"
$ws.Range("C10").Value = "Ultrasonics"
$ws.Range("C11").Value = "Holography"
$ws.Range("C12").Value = "Crystal Structures"
$ws.Range("C13").Value = "Crystal Defects"

# Cells that must display the "x" mark, grouped into contiguous rectangular
# blocks per column so each can be filled with a single Copy/PasteSpecial.
$markRanges = @("D9:D12", "E9:E13", "F9:F13", "G11:G12", "H10:H13", "I9:I12")

foreach ($rng in $markRanges) {
    # Re-use the exact style (font color, alignment, etc.) already used for
    # the "x" marks elsewhere in the sheet (e.g. D3) instead of building a
    # brand new style.
    $ws.Range("D3").Copy()
    $ws.Range($rng).PasteSpecial(-4122)
    $ws.Range($rng).Value = [char]0x00D7
}

# Update the active selection / view to match the saved workbook.
$null = $ws.Range("H13").Select()
